# Add 2022-Q3 data:
#   1. Insert a new "2022-Q3" worksheet (positioned right after "总计",
#      before "2022-Q2") with the quarter's fund-holdings detail.
#   2. Insert a new row at the top of the "总计" (summary) sheet's data
#      table for the 2022-Q3 totals, pushing the older quarters down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" summary sheet - insert the 2022-Q3 summary row at row 2
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()

# Pick up the same look (bold index column, plain data columns) as the
# row immediately below, which still holds the old row-2 formatting.
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 6
$summary.Range("D2").Value = 0.28

# ---------------------------------------------------------------------
# 2) New "2022-Q3" fund-detail sheet, inserted before "2022-Q2"
# ---------------------------------------------------------------------
$anchor = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($anchor)
$q3.Name = "2022-Q3"

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Header row + index column share the workbook's bold/bordered/centered
# look used on every other fund-detail sheet.
$q3.Range("B1:H1").Font.Bold = $true
$q3.Range("B1:H1").HorizontalAlignment = -4108
$q3.Range("B1:H1").VerticalAlignment = -4160
$q3.Range("B1:H1").Borders.LineStyle = 1

$q3.Range("A2:A7").Font.Bold = $true
$q3.Range("A2:A7").HorizontalAlignment = -4108
$q3.Range("A2:A7").VerticalAlignment = -4160
$q3.Range("A2:A7").Borders.LineStyle = 1

# Fund code (B) and the numeric-looking metrics (D:G) are stored as text
# in every existing fund-detail sheet (e.g. leading zeros in fund codes
# must survive) - force text format before writing them.
$q3.Range("B2:B7").NumberFormat = "@"
$q3.Range("D2:G7").NumberFormat = "@"

$data = @(
  @(0, "002601", "中银证券价值精选灵活配置混合", "4.68", "92.76", "4.63", "0.2167", 10),
  @(1, "011269", "中银证券优势制造股票A",       "0.76", "92.92", "4.56", "0.0347", 10),
  @(2, "011270", "中银证券优势制造股票C",       "0.58", "92.92", "4.56", "0.0264", 10),
  @(3, "004250", "银河量化优选混合",             "0.24", "86.30", "1.78", "0.0043", 7),
  @(4, "004913", "中银证券聚瑞混合A",           "0.07", "43.92", "2.28", "0.0016", 9),
  @(5, "004914", "中银证券聚瑞混合C",           "0.05", "43.92", "2.28", "0.0011", 9)
)

$r = 2
foreach ($row in $data) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r++
}
